# Update Ligand/Receptor average & total expression, derived specificity,
# and edge weight/specificity columns with recomputed values from updated TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: RowNumber, then a hashtable of Column -> new value
$updates = @(
    @{ Row = 2; Cells = @{ "G"=2.629231666666667; "H"=7.887695; "I"=0.1414315557047068; "J"=0.1414315557047067; "M"=92.64038833333332; "N"=277.921165; "O"=0.7451295270557885; "P"=0.7451295270557885; "Q"=243.5730426182972; "R"=2192.157383564675; "S"=0.1053848282130126; "T"=0.1053848282130125 } },
    @{ Row = 3; Cells = @{ "G"=2.629231666666667; "H"=7.887695; "I"=0.1414315557047068; "J"=0.1414315557047067; "N"=6.766394999999999; "O"=0.01814126213173672; "P"=0.01814126213173672; "Q"=5.930140001058332; "R"=53.371260009525; "S"=0.002565746925738408; "T"=0.002565746925738408 } },
    @{ Row = 4; Cells = @{ "G"=2.629231666666667; "H"=7.887695; "I"=0.1414315557047068; "J"=0.1414315557047067; "M"=28.72545833333334; "N"=86.17637500000001; "O"=0.2310459570329316; "P"=0.2310459570329316; "Q"=75.52588468951389; "R"=679.7329622056251; "S"=0.03267718914245035; "T"=0.03267718914245034 } },
    @{ Row = 5; Cells = @{ "G"=2.629231666666667; "H"=7.887695; "I"=0.1414315557047068; "J"=0.1414315557047067; "M"=0.706587; "N"=2.119761; "O"=0.005683253779543222; "P"=0.005683253779543223; "Q"=1.857780915655; "R"=16.720028240895; "S"=0.0008037914235054524; "T"=0.0008037914235054523 } },
    @{ Row = 6; Cells = @{ "I"=0.6147160060020365; "J"=0.6147160060020365; "M"=92.64038833333332; "N"=277.921165; "O"=0.7451295270557885; "P"=0.7451295270557885; "Q"=1058.662242538711; "R"=9527.960182848399; "S"=0.4580430468259207; "T"=0.4580430468259207 } },
    @{ Row = 7; Cells = @{ "I"=0.6147160060020365; "J"=0.6147160060020365; "N"=6.766394999999999; "O"=0.01814126213173672; "P"=0.01814126213173672; "Q"=25.77467212546666; "R"=231.9720491292; "S"=0.01115172420145718; "T"=0.01115172420145718 } },
    @{ Row = 8; Cells = @{ "I"=0.6147160060020365; "J"=0.6147160060020365; "M"=28.72545833333334; "N"=86.17637500000001; "O"=0.2310459570329316; "P"=0.2310459570329316; "Q"=328.2645796744445; "R"=2954.38121707; "S"=0.1420276479102018; "T"=0.1420276479102018 } },
    @{ Row = 9; Cells = @{ "I"=0.6147160060020365; "J"=0.6147160060020365; "M"=0.706587; "N"=2.119761; "O"=0.005683253779543222; "P"=0.005683253779543223; "Q"=8.074631285840001; "R"=72.67168157256; "S"=0.003493587064456787; "T"=0.003493587064456788 } },
    @{ Row = 10; Cells = @{ "G"=4.24731; "H"=12.74193; "I"=0.2284711798035388; "J"=0.2284711798035388; "M"=92.64038833333332; "N"=277.921165; "O"=0.7451295270557885; "P"=0.7451295270557885; "Q"=393.4724477720499; "R"=3541.252029948449; "S"=0.1702406221528889; "T"=0.1702406221528889 } },
    @{ Row = 11; Cells = @{ "G"=4.24731; "H"=12.74193; "I"=0.2284711798035388; "J"=0.2284711798035388; "N"=6.766394999999999; "O"=0.01814126213173672; "P"=0.01814126213173672; "Q"=9.579659049149997; "R"=86.21693144234999; "S"=0.004144755562363149; "T"=0.004144755562363149 } },
    @{ Row = 12; Cells = @{ "G"=4.24731; "H"=12.74193; "I"=0.2284711798035388; "J"=0.2284711798035388; "M"=28.72545833333334; "N"=86.17637500000001; "O"=0.2310459570329316; "P"=0.2310459570329316; "Q"=122.00592643375; "R"=1098.05333790375; "S"=0.0527873423921516; "T"=0.05278734239215161 } },
    @{ Row = 13; Cells = @{ "G"=4.24731; "H"=12.74193; "I"=0.2284711798035388; "J"=0.2284711798035388; "M"=0.706587; "N"=2.119761; "O"=0.005683253779543222; "P"=0.005683253779543223; "Q"=3.00109403097; "R"=27.00984627873; "S"=0.001298459696135161; "T"=0.001298459696135161 } },
    @{ Row = 14; Cells = @{ "E"=3; "F"=1; "G"=0.2859396666666667; "H"=0.8578190000000001; "I"=0.01538125848971795; "J"=0.01538125848971795; "M"=92.64038833333332; "N"=277.921165; "O"=0.7451295270557885; "P"=0.7451295270557885; "Q"=26.48956175990389; "R"=238.406055839135; "S"=0.01146102986396637; "T"=0.01146102986396637 } },
    @{ Row = 15; Cells = @{ "E"=3; "F"=1; "G"=0.2859396666666667; "H"=0.8578190000000001; "I"=0.01538125848971795; "J"=0.01538125848971795; "N"=6.766394999999999; "O"=0.01814126213173672; "P"=0.01814126213173672; "Q"=0.6449269102783333; "R"=5.804342192505; "S"=0.0002790354421779742; "T"=0.0002790354421779742 } },
    @{ Row = 16; Cells = @{ "E"=3; "F"=1; "G"=0.2859396666666667; "H"=0.8578190000000001; "I"=0.01538125848971795; "J"=0.01538125848971795; "M"=28.72545833333334; "N"=86.17637500000001; "O"=0.2310459570329316; "P"=0.2310459570329316; "Q"=8.213747980680557; "R"=73.92373182612502; "S"=0.003553777588127788; "T"=0.003553777588127788 } },
    @{ Row = 17; Cells = @{ "E"=3; "F"=1; "G"=0.2859396666666667; "H"=0.8578190000000001; "I"=0.01538125848971795; "J"=0.01538125848971795; "M"=0.706587; "N"=2.119761; "O"=0.005683253779543222; "P"=0.005683253779543223; "Q"=0.202041251251; "R"=1.818371261259; "S"=0.00008741559544582081; "T"=0.00008741559544582083 } }
)

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$row").Value = $update.Cells[$col]
    }
}